# Apply cryptos list update (prices & volume %) scraped on Sun Nov 17 13:55:14 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'90.551.90"
$ws.Range("E2").Value = "'  -0.32%  "

# Row 3
$ws.Range("D3").Value = "'3.102.54"
$ws.Range("E3").Value = "'  -1.93%  "

# Row 4
$ws.Range("E4").Value = "'  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'231.26"
$ws.Range("E5").Value = "'  +5.84%  "

# Row 6
$ws.Range("D6").Value = "'626.87"
$ws.Range("E6").Value = "'  +0.20%  "

# Row 7
$ws.Range("E7").Value = "'  +0.03%  "

# Row 8
$ws.Range("E8").Value = "'  -2.52%  "

# Row 9
$ws.Range("E9").Value = "'  +0.10%  "

# Row 10
$ws.Range("D10").Value = "'3.101.80"
$ws.Range("E10").Value = "'  -1.96%  "

# Row 11
$ws.Range("D11").Value = "'0.726"
$ws.Range("E11").Value = "'  -5.87%  "

# Row 12
$ws.Range("D12").Value = "'0.197"
$ws.Range("E12").Value = "'  -1.54%  "

# Row 13
$ws.Range("D13").Value = "'36.38"
$ws.Range("E13").Value = "'  +2.94%  "

# Row 14
$ws.Range("D14").Value = "'0.0000246"
$ws.Range("E14").Value = "'  -1.31%  "

# Row 15
$ws.Range("D15").Value = "'5.49"
$ws.Range("E15").Value = "'  -2.50%  "

# Row 16
$ws.Range("D16").Value = "'90.347.99"
$ws.Range("E16").Value = "'  -0.35%  "

# Row 17
$ws.Range("D17").Value = "'3.693.33"
$ws.Range("E17").Value = "'  -1.48%  "

# Row 18
$ws.Range("D18").Value = "'3.107.82"
$ws.Range("E18").Value = "'  -0.52%  "

# Row 19
$ws.Range("D19").Value = "'3.81"
$ws.Range("E19").Value = "'  +1.12%  "

# Row 20
$ws.Range("D20").Value = "'14.09"
$ws.Range("E20").Value = "'  -1.57%  "

# Row 21
$ws.Range("D21").Value = "'0.0000210"
$ws.Range("E21").Value = "'  -4.13%  "

# Row 22
$ws.Range("D22").Value = "'440.77"
$ws.Range("E22").Value = "'  -0.40%  "

# Row 23
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "'  +6.81%  "

# Row 24
$ws.Range("D24").Value = "'8.91"
$ws.Range("E24").Value = "'  -0.28%  "

# Row 25
$ws.Range("D25").Value = "'5.88"
$ws.Range("E25").Value = "'  -2.20%  "

# Row 26
$ws.Range("D26").Value = "'89.36"
$ws.Range("E26").Value = "'  +2.47%  "

# Row 27
$ws.Range("D27").Value = "'12.32"
$ws.Range("E27").Value = "'  +0.01%  "

# Row 28
$ws.Range("E28").Value = "'  -0.89%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  -0.05%  "

# Row 30
$ws.Range("D30").Value = "'9.48"
$ws.Range("E30").Value = "'  +1.71%  "

# Row 31
$ws.Range("E31").Value = "'  -1.92%  "

# Row 32
$ws.Range("D32").Value = "'0.204"
$ws.Range("E32").Value = "'  +19.06%  "

# Row 33
$ws.Range("D33").Value = "'26.40"
$ws.Range("E33").Value = "'  +4.59%  "

# Row 34
$ws.Range("E34").Value = "'  -11.85%  "

# Row 35
$ws.Range("D35").Value = "'0.151"
$ws.Range("E35").Value = "'  +4.94%  "

# Row 36
$ws.Range("D36").Value = "'3.81"
$ws.Range("E36").Value = "'  +2.19%  "

# Row 37
$ws.Range("D37").Value = "'509.82"
$ws.Range("E37").Value = "'  -2.66%  "

# Row 38
$ws.Range("E38").Value = "'  +0.60%  "

# Row 39
$ws.Range("D39").Value = "'7.06"
$ws.Range("E39").Value = "'  +0.39%  "

# Row 40
$ws.Range("E40").Value = "'  -1.78%  "

# Row 41
$ws.Range("D41").Value = "'0.0905"
$ws.Range("E41").Value = "'  +7.28%  "

# Row 42
$ws.Range("B42").Value = "'MantraDAO"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").Value = "'3.54"
$ws.Range("E42").Value = "'  +57.27%  "

# Row 43
$ws.Range("B43").Value = "'PolygonEcosystemToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.411"
$ws.Range("E43").Value = "'  -0.01%  "

# Row 44
$ws.Range("B44").Value = "'WhiteBITCoin"
$ws.Range("C44").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'22.18"
$ws.Range("E44").Value = "'  -0.16%  "

# Row 45
$ws.Range("B45").Value = "'USDe"
$ws.Range("C45").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  -0.02%  "

# Row 46
$ws.Range("D46").Value = "'1.91"
$ws.Range("E46").Value = "'  -2.16%  "

# Row 47
$ws.Range("D47").Value = "'151.22"
$ws.Range("E47").Value = "'  +1.32%  "

# Row 48
$ws.Range("D48").Value = "'0.689"
$ws.Range("E48").Value = "'  +6.03%  "

# Row 49
$ws.Range("D49").Value = "'45.08"
$ws.Range("E49").Value = "'  +2.30%  "

# Row 50
$ws.Range("D50").Value = "'1.34"
$ws.Range("E50").Value = "'  -0.12%  "

# Row 51
$ws.Range("D51").Value = "'4.46"
$ws.Range("E51").Value = "'  +1.91%  "
